$d = $word.ActiveDocument

# Locate the end of the sentence that ends with "...on each of our machines.  "
$target = "to have Emeka in our group, who had previously setup an AWS Postgres server account.  Utilizing the correct login credentials, we were able to access this server and forgo a build of the server locally on each of our machines.  "

$range = $d.Content
$found = $range.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Collapse to the end of the found range so we append right after it.
    $range.Collapse(0)
    $range.InsertAfter("There is a transformation step within the flask app:  a join utilizing sqlalchemy in line 48 allows a single dataframe to be searched by the user. ")
}
